$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.3709835751931791
$ws.Range("C2").Value = 0.04126107170033322
$ws.Range("E2").Value = 0.6416241031824086
$ws.Range("F2").Value = 2.342192894683222
$ws.Range("G2").Value = 0.002438789383091667
$ws.Range("J2").Value = 0.050699443948659
$ws.Range("K2").Value = 0.3357123734446077
$ws.Range("N2").Value = 1.41047792437346
$ws.Range("O2").Value = 2.601295450356545
# Row 3
$ws.Range("B3").Value = 0.3344821430482625
$ws.Range("C3").Value = 0.03684668311586847
$ws.Range("E3").Value = 0.6144145472549383
$ws.Range("F3").Value = 2.312097614300939
$ws.Range("G3").Value = 0.002441048068113237
$ws.Range("J3").Value = 0.05111899710836276
$ws.Range("K3").Value = 0.297491681978812
$ws.Range("N3").Value = 1.42571115862351
$ws.Range("O3").Value = 2.62169467694747
# Row 4
$ws.Range("B4").Value = 0.3121277024033304
$ws.Range("C4").Value = 0.03412045392894925
$ws.Range("E4").Value = 0.5980282537050527
$ws.Range("F4").Value = 2.295002006759304
$ws.Range("G4").Value = 0.002442507933518484
$ws.Range("J4").Value = 0.05140657362355405
$ws.Range("K4").Value = 0.2740255471301793
$ws.Range("N4").Value = 1.435577151892058
$ws.Range("O4").Value = 2.635748338849538
# Row 5
$ws.Range("B5").Value = 0.3030330343497667
$ws.Range("C5").Value = 0.0330055764783026
$ws.Range("E5").Value = 0.5914311879223675
$ws.Range("F5").Value = 2.28838276358691
$ws.Range("G5").Value = 0.002443121257906567
$ws.Range("J5").Value = 0.05153129101502785
$ws.Range("K5").Value = 0.2644637527077549
$ws.Range("N5").Value = 1.439726566072082
$ws.Range("O5").Value = 2.641859345981928
# Row 6
$ws.Range("B6").Value = 0.3015237882604822
$ws.Range("C6").Value = 0.03282021683926928
$ws.Range("E6").Value = 0.5903406080554277
$ws.Range("F6").Value = 2.287304610747341
$ws.Range("G6").Value = 0.002443224214214093
$ws.Range("J6").Value = 0.05155245464948877
$ws.Range("K6").Value = 0.2628760887832868
$ws.Range("N6").Value = 1.440423361332908
$ws.Range("O6").Value = 2.642897259060391
# Row 7
$ws.Range("B7").Value = 0.312004987385393
$ws.Range("C7").Value = 0.03410543409444244
$ws.Range("E7").Value = 0.5979389576405936
$ws.Range("F7").Value = 2.294911331455182
$ws.Range("G7").Value = 0.002442516130319012
$ws.Range("J7").Value = 0.05140822513752141
$ws.Range("K7").Value = 0.2738965893873058
$ws.Range("N7").Value = 1.435632590245689
$ws.Range("O7").Value = 2.6358291995477
# Row 8
$ws.Range("B8").Value = 0.3583861981272207
$ws.Range("C8").Value = 0.03974229381073258
$ws.Range("E8").Value = 0.6321756576351163
$ws.Range("F8").Value = 2.331528730364894
$ws.Range("G8").Value = 0.0024395530547347
$ws.Range("J8").Value = 0.05083788070540685
$ws.Range("K8").Value = 0.3225338283303074
$ws.Range("N8").Value = 1.41562393583278
$ws.Range("O8").Value = 2.608011680657654
# Row 9
$ws.Range("B9").Value = 0.4497805965534099
$ws.Range("C9").Value = 0.05066927706027968
$ws.Range("E9").Value = 0.7018656268820251
$ws.Range("F9").Value = 2.414334765379735
$ws.Range("G9").Value = 0.002434319360039165
$ws.Range("J9").Value = 0.04995762492756839
$ws.Range("K9").Value = 0.417907759605157
$ws.Range("N9").Value = 1.380454251689081
$ws.Range("O9").Value = 2.565604576294106
# Row 10
$ws.Range("B10").Value = 0.5171822692200294
$ws.Range("C10").Value = 0.05861844684025641
$ws.Range("E10").Value = 0.7546411100652222
$ws.Range("F10").Value = 2.481923883268934
$ws.Range("G10").Value = 0.002430822306328961
$ws.Range("J10").Value = 0.04945668171342987
$ws.Range("K10").Value = 0.4879627256277672
$ws.Range("N10").Value = 1.357092514772496
$ws.Range("O10").Value = 2.541873606847503
# Row 11
$ws.Range("B11").Value = 0.547897550336387
$ws.Range("C11").Value = 0.06221732337688479
$ws.Range("E11").Value = 0.7789962695346446
$ws.Range("F11").Value = 2.514148561129105
$ws.Range("G11").Value = 0.002429306247082938
$ws.Range("J11").Value = 0.0492605718854513
$ws.Range("K11").Value = 0.5198264511014088
$ws.Range("N11").Value = 1.347002202057642
$ws.Range("O11").Value = 2.532695250060357
# Row 12
$ws.Range("B12").Value = 0.5595360013955997
$ws.Range("C12").Value = 0.06357760466171669
$ws.Range("E12").Value = 0.7882691066418062
$ws.Range("F12").Value = 2.526564409299397
$ws.Range("G12").Value = 0.002428742849451242
$ws.Range("J12").Value = 0.04919088931589499
$ws.Range("K12").Value = 0.5318913752438448
$ws.Range("N12").Value = 1.343258493814762
$ws.Range("O12").Value = 2.529452572508262
# Row 13
$ws.Range("B13").Value = 0.5570291388532098
$ws.Range("C13").Value = 0.06328475738919792
$ws.Range("E13").Value = 0.7862698065650022
$ws.Range("F13").Value = 2.523880949916133
$ws.Range("G13").Value = 0.002428863712102495
$ws.Range("J13").Value = 0.04920569282917242
$ws.Range("K13").Value = 0.529293038447662
$ws.Range("N13").Value = 1.344061330812629
$ws.Range("O13").Value = 2.530140572727902
# Row 14
$ws.Range("B14").Value = 0.5488549105322988
$ws.Range("C14").Value = 0.06232928566485896
$ws.Range("E14").Value = 0.7797581473070636
$ws.Range("F14").Value = 2.515165747352256
$ws.Range("G14").Value = 0.002429259681758042
$ws.Range("J14").Value = 0.04925474720448975
$ws.Range("K14").Value = 0.5208190661900289
$ws.Range("N14").Value = 1.346692655963452
$ws.Range("O14").Value = 2.532423801453945
# Row 15
$ws.Range("B15").Value = 0.5438488873928691
$ws.Range("C15").Value = 0.06174369983742167
$ws.Range("E15").Value = 0.7757760918106982
$ws.Range("F15").Value = 2.50985519662791
$ws.Range("G15").Value = 0.002429503616750511
$ws.Range("J15").Value = 0.04928539120733433
$ws.Range("K15").Value = 0.5156283450379533
$ws.Range("N15").Value = 1.34831448339337
$ws.Range("O15").Value = 2.533852697453824
# Row 16
$ws.Range("B16").Value = 0.5151759945272829
$ws.Range("C16").Value = 0.05838290007856983
$ws.Range("E16").Value = 0.7530564425701272
$ws.Range("F16").Value = 2.479847711283838
$ws.Range("G16").Value = 0.002430922884668263
$ws.Range("J16").Value = 0.04947013823903745
$ws.Range("K16").Value = 0.4858802193305962
$ws.Range("N16").Value = 1.357762751271576
$ws.Range("O16").Value = 2.542506010008992
# Row 17
$ws.Range("B17").Value = 0.4975995488530316
$ws.Range("C17").Value = 0.0563167026403022
$ws.Range("E17").Value = 0.739207718317104
$ws.Range("F17").Value = 2.461818023542463
$ws.Range("G17").Value = 0.002431812673361899
$ws.Range("J17").Value = 0.04959162003392237
$ws.Range("K17").Value = 0.4676291577730751
$ws.Range("N17").Value = 1.363696542726796
$ws.Range("O17").Value = 2.548229002012278
# Row 18
$ws.Range("B18").Value = 0.4874951482255199
$ws.Range("C18").Value = 0.05512666108114672
$ws.Range("E18").Value = 0.7312749831198886
$ws.Range("F18").Value = 2.451586913571077
$ws.Range("G18").Value = 0.002432331496853083
$ws.Range("J18").Value = 0.04966448298173631
$ws.Range("K18").Value = 0.4571312246064281
$ws.Range("N18").Value = 1.367160050163566
$ws.Range("O18").Value = 2.551672886603995
# Row 19
$ws.Range("B19").Value = 0.4840748644093367
$ws.Range("C19").Value = 0.0547234570258297
$ws.Range("E19").Value = 0.7285947052484687
$ws.Range("F19").Value = 2.448146708526963
$ws.Range("G19").Value = 0.002432508372076557
$ws.Range("J19").Value = 0.04968966628954163
$ws.Range("K19").Value = 0.4535767508055244
$ws.Range("N19").Value = 1.368341414520199
$ws.Range("O19").Value = 2.552865048676551
# Row 20
$ws.Range("B20").Value = 0.4994700667629388
$ws.Range("C20").Value = 0.05653682089661061
$ws.Range("E20").Value = 0.7406785550386275
$ws.Range("F20").Value = 2.463722915574948
$ws.Range("G20").Value = 0.002431717225462297
$ws.Range("J20").Value = 0.04957837857180891
$ws.Range("K20").Value = 0.4695720600849427
$ws.Range("N20").Value = 1.363059648691916
$ws.Range("O20").Value = 2.547604027379577
# Row 21
$ws.Range("B21").Value = 0.5512556871133825
$ws.Range("C21").Value = 0.06261000011963347
$ws.Range("E21").Value = 0.7816694196979768
$ws.Range("F21").Value = 2.517719826990799
$ws.Range("G21").Value = 0.002429143086064499
$ws.Range("J21").Value = 0.04924021435665438
$ws.Range("K21").Value = 0.5233081132732309
$ws.Range("N21").Value = 1.345917674401299
$ws.Range("O21").Value = 2.531746835601041
# Row 22
$ws.Range("B22").Value = 0.5851425503008727
$ws.Range("C22").Value = 0.06656437119490022
$ws.Range("E22").Value = 0.8087512049689565
$ws.Range("F22").Value = 2.55425206365922
$ws.Range("G22").Value = 0.002427523088504913
$ws.Range("J22").Value = 0.04904590529783448
$ws.Range("K22").Value = 0.5584206323229921
$ws.Range("N22").Value = 1.335164839997486
$ws.Range("O22").Value = 2.522741362708814
# Row 23
$ws.Range("B23").Value = 0.5670528231289893
$ws.Range("C23").Value = 0.06445522180383989
$ws.Range("E23").Value = 0.7942703981889991
$ws.Range("F23").Value = 2.534640292894608
$ws.Range("G23").Value = 0.002428382023034784
$ws.Range("J23").Value = 0.04914716484699611
$ws.Range("K23").Value = 0.5396812383000054
$ws.Range("N23").Value = 1.340862599103456
$ws.Range("O23").Value = 2.52742333690901
# Row 24
$ws.Range("B24").Value = 0.4986244039197345
$ws.Range("C24").Value = 0.0564373121483186
$ws.Range("E24").Value = 0.7400134991837746
$ws.Range("F24").Value = 2.46286129527752
$ws.Range("G24").Value = 0.002431760354825683
$ws.Range("J24").Value = 0.04958435562572561
$ws.Range("K24").Value = 0.4686936899676084
$ws.Range("N24").Value = 1.363347426221313
$ws.Range("O24").Value = 2.54788609979235
# Row 25
$ws.Range("B25").Value = 0.4250102948025472
$ws.Range("C25").Value = 0.04772698157995592
$ws.Range("E25").Value = 0.6827375343862201
$ws.Range("F25").Value = 2.39075123991968
$ws.Range("G25").Value = 0.00243567382124352
$ws.Range("J25").Value = 0.0501701986288694
$ws.Range("K25").Value = 0.3921084250988827
$ws.Range("N25").Value = 1.389533384779259
